$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Insert a new blank row at row 69 (pushes old rows 69-137 down to 70-138)
# ---------------------------------------------------------------------
$ws.Rows.Item(69).Insert()

# The plain row-insert doesn't fully copy cell borders for the new row in
# this engine, so re-apply formatting from the (untouched) row below,
# which carries the same "normal table data row" style.
$ws.Range("A70:K70").Copy()
$ws.Range("A69:K69").PasteSpecial(-4122)   # xlPasteFormats

# Table1's calculated column ("EARNED ", column G) formula needs to be
# re-asserted on the freshly inserted row.
$ws.Cells.Item(69,7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Grow Table1 so it covers the new trailing row (138) created by the
# downward shift of the former last row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K138"))

# Re-assert the same calculated-column formula on the new last row (138);
# the table resize alone leaves it unresolved.
$ws.Cells.Item(138,7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------
# Data edits
# ---------------------------------------------------------------------
# C67: 1.25 days SL earned, recorded for the 5/2023 period
$ws.Cells.Item(67,3).Value2 = 1.25

# Row 68 (6/2023 period row): SL(1-0-0) taken, 1 day, dated 6/1/2023
$ws.Cells.Item(68,2).Value2 = "SL(1-0-0)"
$ws.Cells.Item(68,8).Value2 = 1
$ws.Cells.Item(67,11).Copy()
$ws.Cells.Item(68,11).PasteSpecial(-4122)  # xlPasteFormats (date number format)
$ws.Cells.Item(68,11).Value2 = 45078

# Row 69 (newly inserted row): another SL(1-0-0) taken, 1 day, dated 6/5/2023
$ws.Cells.Item(69,2).Value2 = "SL(1-0-0)"
$ws.Cells.Item(69,8).Value2 = 1
$ws.Cells.Item(67,11).Copy()
$ws.Cells.Item(69,11).PasteSpecial(-4122)  # xlPasteFormats (date number format)
$ws.Cells.Item(69,11).Value2 = 45082

$excel.CutCopyMode = 0
$excel.CalculateFull()
